$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-Year Summary")

$ws.Range("A3").Value = "Cloud Services"
$ws.Range("A4").Value = "Software Licenses"
$ws.Range("A5").Value = "Support & Maintenance"
$ws.Range("A6").Value = "Professional Services"
$ws.Range("A7").Value = "TOTAL"
